# Commit: "Added new Script RCC0001"
#
# Appends one new test-script row (RCC0001 / OBT / ... / Y) to the bottom
# of the "Test Cases" table on sheet 1, formatted like the existing data
# rows, plus a blank formatted spacer row two rows further down (row 13,
# leaving row 12 empty) - mirroring what Excel leaves behind after a
# fill-down / paste operation that extended the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row (row 11) ---------------------------------------------
$ws.Range("A11").Value = "RCC0001"
$ws.Range("B11").Value = "OBT"
$ws.Range("C11").Value = "Verify that user is able to add an article to the group from search results page."
$ws.Range("D11").Value = "Y"

# B11:C11 take on the same look as the rest of the body rows (B2 carries
# that style - thin bottom border, no fill).
$ws.Range("B2").Copy()
$ws.Range("B11:C11").PasteSpecial(-4122)

# D11 gets a narrower look: only a thin left/right border (no top/bottom).
# Start from the fully-boxed last-row style (D10) and then strip the
# top/bottom edges so just the left/right remain.
$ws.Range("D10").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D11").Borders.Item(9).LineStyle = -4142
$ws.Range("D11").Borders.Item(8).LineStyle = -4142

# --- Blank spacer row (row 13) - same formatting as row 11 -------------
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("B11:D11").Copy()
$ws.Range("B13:D13").PasteSpecial(-4122)

# Leave the selection the way Excel would after entering the spacer row.
$ws.Range("A13:XFD13").Select() | Out-Null
